# Updated the activity status on the "Task Break-Up" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task Break-Up")

# --- Resource (column E) swaps ---------------------------------------
$ws.Range("E5").Value  = "Diwakar"
$ws.Range("E9").Value  = "Sarfaraz"
$ws.Range("E10").Value = "Sarfaraz"
$ws.Range("E11").Value = "Sarfaraz"
$ws.Range("E12").Value = "Sarfaraz"
$ws.Range("E14").Value = "Diwakar"

# --- Status (column G) updates ----------------------------------------
$doneRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,22,35,36,37,38)
foreach ($r in $doneRows) {
    $ws.Range("G$r").Value = "Done"
}

$ws.Range("G15").Value = "In-progress"
$ws.Range("G39").Value = "In-progress"
$ws.Range("G42").Value = "Done"
$ws.Range("G49").Value = "Done "
$ws.Range("G46").Value = "Partialy Done"

# --- Zoom level on this sheet's view ------------------------------------
$ws.Application.ActiveWindow.Zoom = 69
